# Observation Platform vault workbook: add a "types" lookup sheet (CRUD
# support table mapping platform-type codes to English/French labels) in
# front of the existing "Sheet1" data table, and touch up "Sheet1" itself:
# a new "longname" column header, and the "hydrophone" platform type is
# renamed to "mooring".

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "types" worksheet ---------------------------------
# Worksheets.Add() with no arguments inserts the new sheet immediately
# before the active sheet, which is exactly where "types" needs to land
# (ahead of "Sheet1").
$typesSheet = $wb.Worksheets.Add()
$typesSheet.Name = "types"

# Fetch the data sheet reference AFTER inserting "types" so it correctly
# resolves to the "Sheet1" tab (not a stale pre-insert position).
$dataSheet = $wb.Worksheets.Item("Sheet1")

# Column A (numeric ids) for all seven rows.
for ($i = 1; $i -le 7; $i++) {
    $typesSheet.Cells.Item($i, 1).Value = $i
}

# Columns B/C are written in the same order the shared-string table needs
# to see brand-new label text appear (reused labels like "plane"/"boat"/
# "drone" don't matter since they already exist in the string table).
$typesSheet.Cells.Item(1, 2).Value = "plane"
$typesSheet.Cells.Item(1, 3).Value = "avion"

$typesSheet.Cells.Item(2, 2).Value = "boat"
$typesSheet.Cells.Item(2, 3).Value = "bateau"

$typesSheet.Cells.Item(3, 2).Value = "drone"
$typesSheet.Cells.Item(3, 3).Value = "drone"

$typesSheet.Cells.Item(6, 2).Value = "mooring"
$typesSheet.Cells.Item(6, 3).Value = "mouillage"

$typesSheet.Cells.Item(5, 2).Value = "land"
$typesSheet.Cells.Item(5, 3).Value = "terre"

$typesSheet.Cells.Item(7, 2).Value = "space"
$typesSheet.Cells.Item(7, 3).Value = "espace"

$typesSheet.Cells.Item(4, 2).Value = "underwater glider"
$typesSheet.Cells.Item(4, 3).Value = "planeur sous-marin"

[void]$typesSheet.Range("A1:C7").Select()

# --- 2. Update the original "Sheet1" data table ------------------------
# New "longname" header in column F.
$dataSheet.Range("F1").Value = "longname"

# Row 13's platform type was "hydrophone" -- it is now a "mooring".
$dataSheet.Range("A13").Value = "mooring"

[void]$dataSheet.Range("A21").Select()

# --- 3. Make "Sheet1" the active tab -----------------------------------
$dataSheet.Activate()
